# Edit script: apply the changes described by the diff using Word COM interop.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Skills cell: ".NET, MVC, EF, WCF, MS SQL,JavaScript" -> expanded list of
#    runs with new wording/formatting, plus a trailing "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(".NET, MVC, EF, WCF, MS SQL,JavaScript", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r1.Find.Found) {
    throw "Could not find the skills paragraph text"
}
$p1 = $r1.Paragraphs(1).Range
$p1.InsertXML('<w:p w14:paraId="7EF095D8" w14:textId="37425972" w:rsidR="00670AA2" w:rsidRPr="0026750C" w:rsidRDefault="00670AA2" w:rsidP="000D344A"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">.NET, MVC &amp; Web API, EF, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>WCF(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">SOAP &amp; REST), </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">MS </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>SQL, JavaScript</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="333333"/><w:sz w:val="17"/><w:szCs w:val="17"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# ---------------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark (previously sitting in its own empty
#    paragraph after the Extracurricular Activities section).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) "Extracurricular Activities" + ": " runs merge into a single run.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Extracurricular Activities: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r2.Find.Found) {
    throw "Could not find the Extracurricular Activities paragraph text"
}
$p2 = $r2.Paragraphs(1).Range
$p2.InsertXML('<w:p w14:paraId="7020CDDD" w14:textId="0921ACDF" w:rsidR="00FE31F8" w:rsidRPr="00FE31F8" w:rsidRDefault="00FE31F8" w:rsidP="00FE31F8"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00FE31F8"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Extracurricular Activities: </w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 4) Drop the stray <w:lastRenderedPageBreak/> before "Signature".
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Signature", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r3.Find.Found) {
    throw "Could not find the Signature paragraph text"
}
$p3 = $r3.Paragraphs(1).Range
$p3.InsertXML('<w:p w14:paraId="60029506" w14:textId="70C9C80F" w:rsidR="00293B6B" w:rsidRPr="00FC1AF9" w:rsidRDefault="00293B6B" w:rsidP="00571FE9"><w:pPr><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="00947774"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana" w:cs="Calibri"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/></w:rPr><w:t>Signature</w:t></w:r></w:p>')
